$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data-refresh for the "paises" COVID tracking sheet: updated case counts for
# several countries, a couple of countries swap rank position (so their row's
# country-name cell and stats both change), and the "last updated" timestamp
# moves forward. Every touched cell is listed explicitly below; unlisted
# cells are left untouched.
$changes = @(
    @{ Cell = "A1"; Value = 'Datos actualizados a 1 de Octubre de 2020 a las 13:55' },
    @{ Cell = "B4"; Value = 7451354 },
    @{ Cell = "C4"; Value = 4072 },
    @{ Cell = "D4"; Value = 4701240 },
    @{ Cell = "E4"; Value = 2538309 },
    @{ Cell = "G4"; Value = 65 },
    @{ Cell = "H4"; Value = 211805 },
    @{ Cell = "B16"; Value = 461044 },
    @{ Cell = "C16"; Value = 3825 },
    @{ Cell = "D16"; Value = 383368 },
    @{ Cell = "E16"; Value = 51296 },
    @{ Cell = "G16"; Value = 211 },
    @{ Cell = "H16"; Value = 26380 },
    @{ Cell = "B26"; Value = 291182 },
    @{ Cell = "C26"; Value = 4174 },
    @{ Cell = "D26"; Value = 218487 },
    @{ Cell = "E26"; Value = 61839 },
    @{ Cell = "G26"; Value = 116 },
    @{ Cell = "H26"; Value = 10856 },
    @{ Cell = "B40"; Value = 105676 },
    @{ Cell = "C40"; Value = 494 },
    @{ Cell = "D40"; Value = 97197 },
    @{ Cell = "E40"; Value = 7867 },
    @{ Cell = "G40"; Value = 2 },
    @{ Cell = "H40"; Value = 612 },
    @{ Cell = "B43"; Value = 95348 },
    @{ Cell = "C43"; Value = 1158 },
    @{ Cell = "D43"; Value = 84903 },
    @{ Cell = "E43"; Value = 10024 },
    @{ Cell = "G43"; Value = 2 },
    @{ Cell = "H43"; Value = 421 },
    @{ Cell = "A49"; Value = 'Nepal' },
    @{ Cell = "B49"; Value = 79728 },
    @{ Cell = "C49"; Value = 1911 },
    @{ Cell = "D49"; Value = 57389 },
    @{ Cell = "E49"; Value = 21830 },
    @{ Cell = "G49"; Value = 11 },
    @{ Cell = "H49"; Value = 509 },
    @{ Cell = "A50"; Value = 'Bielorrusia' },
    @{ Cell = "B50"; Value = 79019 },
    @{ Cell = "C50"; Value = 388 },
    @{ Cell = "D50"; Value = 74777 },
    @{ Cell = "E50"; Value = 3403 },
    @{ Cell = "G50"; Value = 6 },
    @{ Cell = "H50"; Value = 839 },
    @{ Cell = "B90"; Value = 16454 },
    @{ Cell = "C90"; Value = 46 },
    @{ Cell = "D90"; Value = 15430 },
    @{ Cell = "E90"; Value = 792 },
    @{ Cell = "G90"; Value = 2 },
    @{ Cell = "H90"; Value = 232 },
    @{ Cell = "B91"; Value = 15019 },
    @{ Cell = "C91"; Value = 37 },
    @{ Cell = "D91"; Value = 12538 },
    @{ Cell = "E91"; Value = 2170 },
    @{ Cell = "A100"; Value = 'Eslovaquia' },
    @{ Cell = "B100"; Value = 10938 },
    @{ Cell = "C100"; Value = 797 },
    @{ Cell = "D100"; Value = 4620 },
    @{ Cell = "E100"; Value = 6270 },
    @{ Cell = "H100"; Value = 48 },
    @{ Cell = "A101"; Value = 'Montenegro' },
    @{ Cell = "B101"; Value = 10772 },
    @{ Cell = "D101"; Value = 7192 },
    @{ Cell = "E101"; Value = 3411 },
    @{ Cell = "H101"; Value = 169 },
    @{ Cell = "A102"; Value = 'Consejo Danes para los Refugiados' },
    @{ Cell = "B102"; Value = 10659 },
    @{ Cell = "D102"; Value = 10139 },
    @{ Cell = "E102"; Value = 248 },
    @{ Cell = "H102"; Value = 272 },
    @{ Cell = "A103"; Value = 'Guinea' },
    @{ Cell = "B103"; Value = 10652 },
    @{ Cell = "D103"; Value = 9996 },
    @{ Cell = "E103"; Value = 590 },
    @{ Cell = "H103"; Value = 66 },
    @{ Cell = "A104"; Value = 'Maldivas' },
    @{ Cell = "B104"; Value = 10291 },
    @{ Cell = "D104"; Value = 9108 },
    @{ Cell = "E104"; Value = 1149 },
    @{ Cell = "H104"; Value = 34 },
    @{ Cell = "E115"; Value = 3181 },
    @{ Cell = "G115"; Value = 1 },
    @{ Cell = "H115"; Value = 40 },
    @{ Cell = "E142"; Value = 742 },
    @{ Cell = "G142"; Value = 1 },
    @{ Cell = "H142"; Value = 65 },
    @{ Cell = "B146"; Value = 3095 },
    @{ Cell = "C146"; Value = 37 },
    @{ Cell = "D146"; Value = 2605 },
    @{ Cell = "E146"; Value = 455 },
    @{ Cell = "A207"; Value = 'Santa Lucia' },
    @{ Cell = "A208"; Value = 'Nueva Caledonia' }
)

foreach ($change in $changes) {
    $ws.Range($change.Cell).Value = $change.Value
}
